$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 83 (shifts existing rows 83..122 down to 84..123,
# growing the used range from A1:T122 to A1:T123), then populate the new row with
# the new weekly price observation.
$ws.Rows.Item(83).Insert()

$ws.Cells.Item(83, 1).Value = 1
$ws.Cells.Item(83, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(83, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(83, 4).Value = 44839
$ws.Cells.Item(83, 5).Value = 15
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100102
$ws.Cells.Item(83, 8).Value = "Cítricos"
$ws.Cells.Item(83, 9).Value = 100102004
$ws.Cells.Item(83, 10).Value = "Mandarina"
$ws.Cells.Item(83, 11).Value = "Murcott"
$ws.Cells.Item(83, 12).Value = "Segunda"
$ws.Cells.Item(83, 13).Value = 250
$ws.Cells.Item(83, 14).Value = 16000
$ws.Cells.Item(83, 15).Value = 17000
$ws.Cells.Item(83, 16).Value = 16500
$ws.Cells.Item(83, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(83, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(83, 19).Value = 825
$ws.Cells.Item(83, 20).Value = 20
